# Update the "Förändrad" (Changed) date column (C) for all data rows
# from the old serial date value 45182 (2023-09-13) to the new value
# 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Data rows run from row 2 to row 262 (row 1 is the header row).
$ws.Range("C2:C262").Value = 45184
